# Update the "dSF" column (F) values for several rows in Sheet1.
# This reflects a repull/recalculation of the dSF data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -3
$ws.Range("F6").Value = -6
$ws.Range("F8").Value = -9
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -3
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = -6
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -1
